# Remove the stray "Market Failure" placeholder slide (sldId 267) that was
# accidentally added as slide 2, and remove the picture that was accidentally
# added to the title slide (slide 1).

$p = $ppt.ActivePresentation

# Delete the picture shape that was added to the title slide (slide 1).
$titleSlide = $p.Slides.Item(1)
for ($i = $titleSlide.Shapes.Count; $i -ge 1; $i--) {
    $shape = $titleSlide.Shapes.Item($i)
    if ($shape.Type -eq 13) {
        # msoPicture
        $shape.Delete()
    }
}

# Delete the extra blank slide (originally slide 2, SlideID 267) that only
# held the now-removed picture.
for ($i = $p.Slides.Count; $i -ge 1; $i--) {
    $slide = $p.Slides.Item($i)
    if ($slide.SlideID -eq 267) {
        $slide.Delete()
    }
}
